$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B3").Value = 'Pabellón De Arteaga'
$ws.Range("B4").Value = 'Rincón De Romos'
$ws.Range("B21").Value = 'Chiapa De Corzo'
$ws.Range("B40").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B44").Value = 'San Cristóbal De Las Casas'
$ws.Range("A76").Value = 'Ciudad De México'
$ws.Range("B80").Value = 'Cuajimalpa De Morelos'
$ws.Range("D81").Value = 0.009102730819245772
$ws.Range("B99").Value = 'San Juan Del Río'
$ws.Range("A103").Value = 'Estado De México'
$ws.Range("B103").Value = 'Almoloya De Alquisiras'
$ws.Range("B104").Value = 'Almoloya De Juárez'
$ws.Range("B108").Value = 'Atizapán De Zaragoza'
$ws.Range("B114").Value = 'Chapa De Mota'
$ws.Range("B120").Value = 'Ecatepec De Morelos'
$ws.Range("B130").Value = 'Naucalpan De Juárez'
$ws.Range("B136").Value = 'San Felipe Del Progreso'
$ws.Range("B146").Value = 'Tlalnepantla De Baz'
$ws.Range("B150").Value = 'Valle De Bravo'
$ws.Range("B151").Value = 'Villa De Allende'
$ws.Range("B152").Value = 'Villa Del Carbón'
$ws.Range("B157").Value = 'Apaseo El Alto'
$ws.Range("B162").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B165").Value = 'Jaral Del Progreso'
$ws.Range("B173").Value = 'San Diego De La Unión'
$ws.Range("B176").Value = 'San Luis De La Paz'
$ws.Range("B177").Value = 'Valle De Santiago'
$ws.Range("B181").Value = 'Acapulco De Juárez'
$ws.Range("B185").Value = 'Atenango Del Río'
$ws.Range("B186").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B188").Value = 'Atoyac De Álvarez'
$ws.Range("B189").Value = 'Ayutla De Los Libres'
$ws.Range("B191").Value = 'Chilapa De Álvarez'
$ws.Range("B192").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B193").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B196").Value = 'Coyuca De Benítez'
$ws.Range("B197").Value = 'Coyuca De Catalán'
$ws.Range("B200").Value = 'Cuetzala Del Progreso'
$ws.Range("B201").Value = 'Cutzamala De Pinzón'
$ws.Range("B205").Value = 'Iguala De La Independencia'
$ws.Range("B208").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B220").Value = 'Taxco De Alarcón'
$ws.Range("B222").Value = 'Técpan De Galeana'
$ws.Range("B224").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B226").Value = 'Tixtla De Guerrero'
$ws.Range("B229").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B230").Value = 'Tlapa De Comonfort'
$ws.Range("B242").Value = 'Atotonilco El Grande'
$ws.Range("B247").Value = 'Cuautepec De Hinojosa'
$ws.Range("B250").Value = 'Huasca De Ocampo'
$ws.Range("B252").Value = 'Huejutla De Reyes'
$ws.Range("B255").Value = 'Jacala De Ledezma'
$ws.Range("B260").Value = 'Omitlán De Juárez'
$ws.Range("B261").Value = 'Pachuca De Soto'
$ws.Range("B266").Value = 'Tepehuacán De Guerrero'
$ws.Range("B268").Value = 'Tezontepec De Aldama'
$ws.Range("B271").Value = 'Tula De Allende'
$ws.Range("B272").Value = 'Tulancingo De Bravo'
$ws.Range("B274").Value = 'Zacualtipán De Ángeles'
$ws.Range("B277").Value = 'Autlán De Navarro'
$ws.Range("B279").Value = 'Encarnación De Díaz'
$ws.Range("B281").Value = 'Lagos De Moreno'
$ws.Range("B284").Value = 'San Cristóbal De La Barranca'
$ws.Range("B288").Value = 'Tizapán El Alto'
$ws.Range("B291").Value = 'Unión De San Antonio'
$ws.Range("B294").Value = 'Zapotlán El Grande'
$ws.Range("B330").Value = 'Coatlán Del Río'
$ws.Range("B336").Value = 'Puente De Ixtla'
$ws.Range("B338").Value = 'Tetela Del Volcán'
$ws.Range("B343").Value = 'Zacualpan De Amilpas'
$ws.Range("B345").Value = 'Ixtlán Del Río'
$ws.Range("B362").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B365").Value = 'Fresnillo De Trujano'
$ws.Range("B367").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B368").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B369").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B371").Value = 'Ixtlán De Juárez'
$ws.Range("B372").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B375").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B376").Value = 'Nejapa De Madero'
$ws.Range("B377").Value = 'Oaxaca De Juárez'
$ws.Range("B378").Value = 'Ocotlán De Morelos'
$ws.Range("B379").Value = 'Putla Villa De Guerrero'
$ws.Range("B396").Value = 'San José Del Progreso'
$ws.Range("B415").Value = 'San Martín De Los Cansecos'
$ws.Range("B432").Value = 'Santa Ana Del Valle'
$ws.Range("B436").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B443").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B463").Value = 'Teotitlán De Flores Magón'
$ws.Range("B464").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B465").Value = 'Tlacolula De Matamoros'
$ws.Range("B467").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B468").Value = 'Villa De Zaachila'
$ws.Range("B469").Value = 'Villa Sola De Vega'
$ws.Range("B479").Value = 'Ayotoxco De Guerrero'
$ws.Range("B480").Value = 'Chalchicomula De Sesma'
$ws.Range("B487").Value = 'Chila De La Sal'
$ws.Range("B491").Value = 'Cuetzalan Del Progreso'
$ws.Range("B498").Value = 'Huehuetlán El Chico'
$ws.Range("B501").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B504").Value = 'Izúcar De Matamoros'
$ws.Range("B512").Value = 'Los Reyes De Juárez'
$ws.Range("B516").Value = 'Palmar De Bravo'
$ws.Range("B529").Value = 'San Salvador El Seco'
$ws.Range("B533").Value = 'Tecali De Herrera'
$ws.Range("B538").Value = 'Tepanco De López'
$ws.Range("B539").Value = 'Tepango De Rodríguez'
$ws.Range("B540").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B545").Value = 'Tepexi De Rodríguez'
$ws.Range("B548").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B567").Value = 'Cadereyta De Montes'
$ws.Range("B571").Value = 'Jalpan De Serra'
$ws.Range("B572").Value = 'Pinal De Amoles'
$ws.Range("B580").Value = 'Ciudad Del Maíz'
$ws.Range("B594").Value = 'Tanquián De Escobedo'
$ws.Range("B595").Value = 'Villa De Arista'
$ws.Range("B596").Value = 'Villa De Arriaga'
$ws.Range("B597").Value = 'Villa De Guadalupe'
$ws.Range("B623").Value = 'Soto La Marina'
$ws.Range("B629").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B634").Value = 'San Pablo Del Monte'
$ws.Range("B639").Value = 'Tetla De La Solidaridad'
$ws.Range("B662").Value = 'Cosamaloapan De Carpio'
$ws.Range("B670").Value = 'Hueyapan De Ocampo'
$ws.Range("B671").Value = 'Ignacio De La Llave'
$ws.Range("B674").Value = 'Ixhuatlán De Madero'
$ws.Range("B675").Value = 'Ixhuatlán Del Café'
$ws.Range("B678").Value = 'Juchique De Ferrer'
$ws.Range("B681").Value = 'Lerdo De Tejada'
$ws.Range("B682").Value = 'Martínez De La Torre'
$ws.Range("B684").Value = 'Medellín De Bravo'
$ws.Range("B691").Value = 'Paso De Ovejas'
$ws.Range("B693").Value = 'Poza Rica De Hidalgo'
$ws.Range("B697").Value = 'Sayula De Alemán'
$ws.Range("D718").Value = 0.09492847854356308
$ws.Range("B729").Value = 'Noria De Ángeles'

# Remove footer/metadata rows 742-746
$ws.Rows("742:746").Delete()
